$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2430.628
$ws.Cells.Item(15, 9).Value = 2430.628
$ws.Cells.Item(15, 11).Value = 7291.884
$ws.Cells.Item(15, 13).Value = -7122.884
$ws.Cells.Item(18, 8).Value = 1029.5555
$ws.Cells.Item(18, 9).Value = 658.25
$ws.Cells.Item(18, 11).Value = 658.25
$ws.Cells.Item(18, 13).Value = -374.25
$ws.Cells.Item(32, 8).Value = 6947.067
$ws.Cells.Item(32, 9).Value = 12360
$ws.Cells.Item(32, 10).Value = 6114.3076
$ws.Cells.Item(32, 11).Value = 12360
$ws.Cells.Item(32, 12).Value = 6114.3076
$ws.Cells.Item(32, 13).Value = -12034
$ws.Cells.Item(32, 14).Value = -6766.3076
$ws.Cells.Item(40, 8).Value = 9526.130999999999
$ws.Cells.Item(40, 9).Value = 5625.1665
$ws.Cells.Item(40, 11).Value = 5625.1665
$ws.Cells.Item(40, 13).Value = -5450.1665
$ws.Cells.Item(43, 8).Value = 10442.25
$ws.Cells.Item(43, 9).Value = 3250
$ws.Cells.Item(43, 10).Value = 11880.7
$ws.Cells.Item(43, 11).Value = 3250
$ws.Cells.Item(43, 12).Value = 11880.7
$ws.Cells.Item(43, 13).Value = -3181
$ws.Cells.Item(43, 14).Value = -12018.7
$ws.Cells.Item(45, 8).Value = 6205
$ws.Cells.Item(45, 9).Value = 475
$ws.Cells.Item(45, 10).Value = 8497
$ws.Cells.Item(45, 11).Value = 1425
$ws.Cells.Item(45, 12).Value = 25491
$ws.Cells.Item(45, 13).Value = -1233
$ws.Cells.Item(45, 14).Value = -25875
$ws.Cells.Item(51, 8).Value = 8947.166999999999
$ws.Cells.Item(51, 10).Value = 9339.933999999999
$ws.Cells.Item(51, 12).Value = 9339.933999999999
$ws.Cells.Item(51, 14).Value = -10307.934
$ws.Cells.Item(62, 8).Value = 8077.3076
$ws.Cells.Item(62, 9).Value = 5250
$ws.Cells.Item(62, 11).Value = 5250
$ws.Cells.Item(62, 13).Value = -4626
$ws.Cells.Item(65, 8).Value = 8077.3076
$ws.Cells.Item(65, 9).Value = 5250
$ws.Cells.Item(65, 11).Value = 26250
$ws.Cells.Item(65, 13).Value = -23130
$ws.Cells.Item(80, 8).Value = 3291.5
$ws.Cells.Item(80, 9).Value = 2499
$ws.Cells.Item(80, 10).Value = 3450
$ws.Cells.Item(80, 11).Value = 7497
$ws.Cells.Item(80, 12).Value = 10350
$ws.Cells.Item(80, 13).Value = -6499
$ws.Cells.Item(80, 14).Value = -12346
$ws.Cells.Item(83, 8).Value = 3291.5
$ws.Cells.Item(83, 9).Value = 2499
$ws.Cells.Item(83, 10).Value = 3450
$ws.Cells.Item(83, 11).Value = 22491
$ws.Cells.Item(83, 12).Value = 31050
$ws.Cells.Item(83, 13).Value = -17499
$ws.Cells.Item(83, 14).Value = -41034
$ws.Cells.Item(98, 8).Value = 20850.7
$ws.Cells.Item(98, 9).Value = 1214.4286
$ws.Cells.Item(98, 10).Value = 66668.664
$ws.Cells.Item(98, 11).Value = 1214.4286
$ws.Cells.Item(98, 12).Value = 66668.664
$ws.Cells.Item(98, 13).Value = 283.5714
$ws.Cells.Item(98, 14).Value = -69664.664
$ws.Cells.Item(116, 8).Value = 9326.462
$ws.Cells.Item(116, 9).Value = 6791.5
$ws.Cells.Item(116, 10).Value = 11499.286
$ws.Cells.Item(116, 11).Value = 6791.5
$ws.Cells.Item(116, 12).Value = 11499.286
$ws.Cells.Item(116, 13).Value = -3349.5
$ws.Cells.Item(116, 14).Value = -18383.286
$ws.Cells.Item(122, 8).Value = 20850.7
$ws.Cells.Item(122, 9).Value = 1214.4286
$ws.Cells.Item(122, 10).Value = 66668.664
$ws.Cells.Item(122, 11).Value = 3643.2858
$ws.Cells.Item(122, 12).Value = 200005.992
$ws.Cells.Item(122, 13).Value = -1193.2858
$ws.Cells.Item(122, 14).Value = -204905.992
$ws.Cells.Item(127, 8).Value = 907.9091
$ws.Cells.Item(127, 9).Value = 548.7
$ws.Cells.Item(127, 11).Value = 1646.1
$ws.Cells.Item(127, 13).Value = 3313.9
$ws.Cells.Item(131, 8).Value = 8878.299999999999
$ws.Cells.Item(131, 10).Value = 6266.3335
$ws.Cells.Item(131, 12).Value = 18799.0005
$ws.Cells.Item(131, 14).Value = -28879.0005
$ws.Cells.Item(132, 8).Value = 3433.25
$ws.Cells.Item(132, 9).Value = 3452.087
$ws.Cells.Item(132, 11).Value = 10356.261
$ws.Cells.Item(132, 13).Value = -7826.261
$ws.Cells.Item(136, 8).Value = 68989.35000000001
$ws.Cells.Item(136, 10).Value = 68989.35000000001
$ws.Cells.Item(136, 12).Value = 68989.35000000001
$ws.Cells.Item(136, 14).Value = -79189.35000000001
$ws.Cells.Item(137, 8).Value = 3097.875
$ws.Cells.Item(137, 9).Value = 2184.4443
$ws.Cells.Item(137, 10).Value = 3363.0645
$ws.Cells.Item(137, 11).Value = 6553.3329
$ws.Cells.Item(137, 12).Value = 10089.1935
$ws.Cells.Item(137, 13).Value = -4003.3329
$ws.Cells.Item(137, 14).Value = -15189.1935
$ws.Cells.Item(138, 8).Value = 2995.7188
$ws.Cells.Item(138, 9).Value = 2252.5
$ws.Cells.Item(138, 10).Value = 4630.8
$ws.Cells.Item(138, 11).Value = 6757.5
$ws.Cells.Item(138, 12).Value = 13892.4
$ws.Cells.Item(138, 13).Value = -1617.5
$ws.Cells.Item(138, 14).Value = -24172.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1745.807
$ws.Cells.Item(32, 9).Value = 1520.1177
$ws.Cells.Item(32, 11).Value = 1520.1177
$ws.Cells.Item(32, 13).Value = -1233.1177
$ws.Cells.Item(61, 8).Value = 7313.143
$ws.Cells.Item(61, 9).Value = 6192.5625
$ws.Cells.Item(61, 10).Value = 8807.25
$ws.Cells.Item(61, 11).Value = 6192.5625
$ws.Cells.Item(61, 12).Value = 8807.25
$ws.Cells.Item(61, 13).Value = -5980.5625
$ws.Cells.Item(61, 14).Value = -9231.25
$ws.Cells.Item(74, 8).Value = 8335769
$ws.Cells.Item(74, 9).Value = 8335769
$ws.Cells.Item(74, 11).Value = 8335769
$ws.Cells.Item(74, 13).Value = -8334895
$ws.Cells.Item(77, 8).Value = 8335769
$ws.Cells.Item(77, 9).Value = 8335769
$ws.Cells.Item(77, 11).Value = 41678845
$ws.Cells.Item(77, 13).Value = -41674477
$ws.Cells.Item(97, 8).Value = 621.8205
$ws.Cells.Item(97, 9).Value = 568.9459000000001
$ws.Cells.Item(97, 11).Value = 568.9459000000001
$ws.Cells.Item(97, 13).Value = -72.94590000000005
$ws.Cells.Item(102, 8).Value = 2821.75
$ws.Cells.Item(102, 9).Value = 2796.2856
$ws.Cells.Item(102, 11).Value = 2796.2856
$ws.Cells.Item(102, 13).Value = -1174.2856
$ws.Cells.Item(122, 8).Value = 4469.5654
$ws.Cells.Item(122, 9).Value = 3531.077
$ws.Cells.Item(122, 11).Value = 10593.231
$ws.Cells.Item(122, 13).Value = -8143.231
$ws.Cells.Item(132, 8).Value = 1794.721
$ws.Cells.Item(132, 9).Value = 924.375
$ws.Cells.Item(132, 10).Value = 4326.636
$ws.Cells.Item(132, 11).Value = 2773.125
$ws.Cells.Item(132, 12).Value = 12979.908
$ws.Cells.Item(132, 13).Value = -243.125
$ws.Cells.Item(132, 14).Value = -18039.908
$ws.Cells.Item(136, 8).Value = 7313.143
$ws.Cells.Item(136, 9).Value = 6192.5625
$ws.Cells.Item(136, 10).Value = 8807.25
$ws.Cells.Item(136, 11).Value = 18577.6875
$ws.Cells.Item(136, 12).Value = 26421.75
$ws.Cells.Item(136, 13).Value = -16027.6875
$ws.Cells.Item(136, 14).Value = -31521.75
$ws.Cells.Item(139, 8).Value = 83160.664
$ws.Cells.Item(139, 10).Value = 83160.664
$ws.Cells.Item(139, 12).Value = 83160.664
$ws.Cells.Item(139, 14).Value = -93440.664
$ws.Cells.Item(140, 8).Value = 88550.336
$ws.Cells.Item(140, 10).Value = 88550.336
$ws.Cells.Item(140, 12).Value = 88550.336
$ws.Cells.Item(140, 14).Value = -98910.336

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 1466.6666
$ws.Cells.Item(7, 9).Value = 200
$ws.Cells.Item(7, 10).Value = 2100
$ws.Cells.Item(7, 11).Value = 200
$ws.Cells.Item(7, 12).Value = 2100
$ws.Cells.Item(7, 13).Value = -87
$ws.Cells.Item(7, 14).Value = -2326
$ws.Cells.Item(20, 8).Value = 2092
$ws.Cells.Item(20, 9).Value = 1100.3334
$ws.Cells.Item(20, 11).Value = 1100.3334
$ws.Cells.Item(20, 13).Value = -853.3334
$ws.Cells.Item(86, 8).Value = 4735.25
$ws.Cells.Item(86, 9).Value = 2495
$ws.Cells.Item(86, 10).Value = 8469
$ws.Cells.Item(86, 11).Value = 2495
$ws.Cells.Item(86, 12).Value = 8469
$ws.Cells.Item(86, 13).Value = -1372
$ws.Cells.Item(86, 14).Value = -10715
$ws.Cells.Item(89, 8).Value = 4735.25
$ws.Cells.Item(89, 9).Value = 2495
$ws.Cells.Item(89, 10).Value = 8469
$ws.Cells.Item(89, 11).Value = 12475
$ws.Cells.Item(89, 12).Value = 42345
$ws.Cells.Item(89, 13).Value = -6859
$ws.Cells.Item(89, 14).Value = -53577
$ws.Cells.Item(105, 8).Value = 26693.555
$ws.Cells.Item(105, 9).Value = 29207.857
$ws.Cells.Item(105, 11).Value = 29207.857
$ws.Cells.Item(105, 13).Value = -27460.857
$ws.Cells.Item(134, 8).Value = 3945
$ws.Cells.Item(134, 9).Value = 1710.6957
$ws.Cells.Item(134, 10).Value = 9654.888999999999
$ws.Cells.Item(134, 11).Value = 5132.0871
$ws.Cells.Item(134, 12).Value = 28964.667
$ws.Cells.Item(134, 13).Value = -2597.0871
$ws.Cells.Item(134, 14).Value = -34034.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 21804.56
$ws.Cells.Item(31, 9).Value = 3387.7058
$ws.Cells.Item(31, 10).Value = 29631.725
$ws.Cells.Item(31, 11).Value = 3387.7058
$ws.Cells.Item(31, 12).Value = 29631.725
$ws.Cells.Item(31, 13).Value = -3092.7058
$ws.Cells.Item(31, 14).Value = -30221.725
$ws.Cells.Item(34, 8).Value = 21804.56
$ws.Cells.Item(34, 9).Value = 3387.7058
$ws.Cells.Item(34, 10).Value = 29631.725
$ws.Cells.Item(34, 11).Value = 3387.7058
$ws.Cells.Item(34, 12).Value = 29631.725
$ws.Cells.Item(34, 13).Value = -3185.7058
$ws.Cells.Item(34, 14).Value = -30035.725
$ws.Cells.Item(58, 8).Value = 6278.95
$ws.Cells.Item(58, 9).Value = 4330.5
$ws.Cells.Item(58, 10).Value = 8227.4
$ws.Cells.Item(58, 11).Value = 4330.5
$ws.Cells.Item(58, 12).Value = 8227.4
$ws.Cells.Item(58, 13).Value = -4127.5
$ws.Cells.Item(58, 14).Value = -8633.4
$ws.Cells.Item(122, 8).Value = 9763.883
$ws.Cells.Item(122, 9).Value = 5229.875
$ws.Cells.Item(122, 10).Value = 13794.111
$ws.Cells.Item(122, 11).Value = 15689.625
$ws.Cells.Item(122, 12).Value = 41382.333
$ws.Cells.Item(122, 13).Value = -13239.625
$ws.Cells.Item(122, 14).Value = -46282.333
$ws.Cells.Item(136, 8).Value = 6278.95
$ws.Cells.Item(136, 9).Value = 4330.5
$ws.Cells.Item(136, 10).Value = 8227.4
$ws.Cells.Item(136, 11).Value = 12991.5
$ws.Cells.Item(136, 12).Value = 24682.2
$ws.Cells.Item(136, 13).Value = -10441.5
$ws.Cells.Item(136, 14).Value = -29782.2
$ws.Cells.Item(140, 8).Value = 116663.836
$ws.Cells.Item(140, 10).Value = 116663.836
$ws.Cells.Item(140, 12).Value = 116663.836
$ws.Cells.Item(140, 14).Value = -127023.836

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 178019.89
$ws.Cells.Item(2, 9).Value = 225.4
$ws.Cells.Item(2, 11).Value = 1352.4
$ws.Cells.Item(2, 13).Value = -1239.4
$ws.Cells.Item(4, 8).Value = 4619483.5
$ws.Cells.Item(4, 9).Value = 3511756.8
$ws.Cells.Item(4, 10).Value = 6219533.5
$ws.Cells.Item(4, 11).Value = 10535270.4
$ws.Cells.Item(4, 12).Value = 18658600.5
$ws.Cells.Item(4, 13).Value = -10535158.4
$ws.Cells.Item(4, 14).Value = -18658824.5
$ws.Cells.Item(68, 8).Value = 2956.963
$ws.Cells.Item(68, 9).Value = 1019.4
$ws.Cells.Item(68, 10).Value = 3397.318
$ws.Cells.Item(68, 11).Value = 3058.2
$ws.Cells.Item(68, 12).Value = 10191.954
$ws.Cells.Item(68, 13).Value = -2247.2
$ws.Cells.Item(68, 14).Value = -11813.954
$ws.Cells.Item(71, 8).Value = 2956.963
$ws.Cells.Item(71, 9).Value = 1019.4
$ws.Cells.Item(71, 10).Value = 3397.318
$ws.Cells.Item(71, 11).Value = 9174.6
$ws.Cells.Item(71, 12).Value = 30575.862
$ws.Cells.Item(71, 13).Value = -5118.6
$ws.Cells.Item(71, 14).Value = -38687.862
$ws.Cells.Item(107, 8).Value = 1772.8636
$ws.Cells.Item(107, 10).Value = 2289.3572
$ws.Cells.Item(107, 12).Value = 6868.071599999999
$ws.Cells.Item(107, 14).Value = -10708.0716
$ws.Cells.Item(108, 8).Value = 6756.4546
$ws.Cells.Item(108, 9).Value = 2375.1667
$ws.Cells.Item(108, 11).Value = 7125.500100000001
$ws.Cells.Item(108, 13).Value = -4245.500100000001
$ws.Cells.Item(110, 8).Value = 10933
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 8105259
$ws.Cells.Item(131, 9).Value = 4478
$ws.Cells.Item(131, 10).Value = 10237043
$ws.Cells.Item(131, 11).Value = 13434
$ws.Cells.Item(131, 12).Value = 30711129
$ws.Cells.Item(131, 13).Value = -8394
$ws.Cells.Item(131, 14).Value = -30721209

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 6
$ws.Cells.Item(21, 9).Value = 6
$ws.Cells.Item(21, 11).Value = 6
$ws.Cells.Item(21, 13).Value = 167
$ws.Cells.Item(24, 8).Value = 37604.11
$ws.Cells.Item(24, 9).Value = 28800.8
$ws.Cells.Item(24, 11).Value = 28800.8
$ws.Cells.Item(24, 13).Value = -28627.8
$ws.Cells.Item(30, 8).Value = 6
$ws.Cells.Item(30, 9).Value = 6
$ws.Cells.Item(30, 11).Value = 6
$ws.Cells.Item(30, 13).Value = 99
$ws.Cells.Item(70, 8).Value = 6013.8823
$ws.Cells.Item(70, 9).Value = 5641.9
$ws.Cells.Item(70, 11).Value = 5641.9
$ws.Cells.Item(70, 13).Value = -5371.9
$ws.Cells.Item(73, 8).Value = 6013.8823
$ws.Cells.Item(73, 9).Value = 5641.9
$ws.Cells.Item(73, 11).Value = 5641.9
$ws.Cells.Item(73, 13).Value = -4705.9
$ws.Cells.Item(102, 8).Value = 2522.0908
$ws.Cells.Item(102, 9).Value = 1814.0714
$ws.Cells.Item(102, 11).Value = 1814.0714
$ws.Cells.Item(102, 13).Value = -192.0714
$ws.Cells.Item(122, 8).Value = 6578.303
$ws.Cells.Item(122, 9).Value = 3930.611
$ws.Cells.Item(122, 10).Value = 9755.532999999999
$ws.Cells.Item(122, 11).Value = 11791.833
$ws.Cells.Item(122, 12).Value = 29266.599
$ws.Cells.Item(122, 13).Value = -9341.832999999999
$ws.Cells.Item(122, 14).Value = -34166.599
$ws.Cells.Item(132, 8).Value = 8968.130999999999
$ws.Cells.Item(132, 9).Value = 1754.6
$ws.Cells.Item(132, 10).Value = 10971.889
$ws.Cells.Item(132, 11).Value = 5263.799999999999
$ws.Cells.Item(132, 12).Value = 32915.667
$ws.Cells.Item(132, 13).Value = -2733.799999999999
$ws.Cells.Item(132, 14).Value = -37975.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 8600.125
$ws.Cells.Item(22, 9).Value = 1225
$ws.Cells.Item(22, 10).Value = 11058.5
$ws.Cells.Item(22, 11).Value = 1225
$ws.Cells.Item(22, 12).Value = 11058.5
$ws.Cells.Item(22, 13).Value = -930
$ws.Cells.Item(22, 14).Value = -11648.5
$ws.Cells.Item(27, 8).Value = 8600.125
$ws.Cells.Item(27, 9).Value = 1225
$ws.Cells.Item(27, 10).Value = 11058.5
$ws.Cells.Item(27, 11).Value = 1225
$ws.Cells.Item(27, 12).Value = 11058.5
$ws.Cells.Item(27, 13).Value = -1118
$ws.Cells.Item(27, 14).Value = -11272.5
$ws.Cells.Item(40, 8).Value = 11666.833
$ws.Cells.Item(40, 9).Value = 6332.3335
$ws.Cells.Item(40, 11).Value = 6332.3335
$ws.Cells.Item(40, 13).Value = -6196.3335
$ws.Cells.Item(46, 8).Value = 4446
$ws.Cells.Item(46, 10).Value = 5319.8
$ws.Cells.Item(46, 12).Value = 5319.8
$ws.Cells.Item(46, 14).Value = -5695.8
$ws.Cells.Item(122, 8).Value = 7453.375
$ws.Cells.Item(122, 9).Value = 5804.091
$ws.Cells.Item(122, 11).Value = 17412.273
$ws.Cells.Item(122, 13).Value = -14962.273
$ws.Cells.Item(132, 8).Value = 4289.7085
$ws.Cells.Item(132, 9).Value = 3111.2727
$ws.Cells.Item(132, 10).Value = 17252.5
$ws.Cells.Item(132, 11).Value = 9333.8181
$ws.Cells.Item(132, 12).Value = 51757.5
$ws.Cells.Item(132, 13).Value = -6803.8181
$ws.Cells.Item(132, 14).Value = -56817.5
$ws.Cells.Item(136, 8).Value = 6896.2705
$ws.Cells.Item(136, 9).Value = 4817.5186
$ws.Cells.Item(136, 10).Value = 12508.9
$ws.Cells.Item(136, 11).Value = 14452.5558
$ws.Cells.Item(136, 12).Value = 37526.7
$ws.Cells.Item(136, 13).Value = -11902.5558
$ws.Cells.Item(136, 14).Value = -42626.7
$ws.Cells.Item(137, 8).Value = 68332.664
$ws.Cells.Item(137, 9).Value = 65000
$ws.Cells.Item(137, 10).Value = 69999
$ws.Cells.Item(137, 11).Value = 65000
$ws.Cells.Item(137, 12).Value = 69999
$ws.Cells.Item(137, 13).Value = -59900
$ws.Cells.Item(137, 14).Value = -80199
$ws.Cells.Item(139, 8).Value = 69998.5
$ws.Cells.Item(139, 10).Value = 69998.5
$ws.Cells.Item(139, 12).Value = 69998.5
$ws.Cells.Item(139, 14).Value = -80278.5
$ws.Cells.Item(141, 8).Value = 69999
$ws.Cells.Item(141, 10).Value = 69999
$ws.Cells.Item(141, 12).Value = 69999
$ws.Cells.Item(141, 14).Value = -80359

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 38114
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 38114
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 38114
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -38610
$ws.Cells.Item(67, 8).Value = 38114
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 38114
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 38114
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -39830
$ws.Cells.Item(81, 8).Value = 5781.8184
$ws.Cells.Item(81, 9).Value = 2325
$ws.Cells.Item(81, 10).Value = 15000
$ws.Cells.Item(81, 11).Value = 4650
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 13).Value = -3589
$ws.Cells.Item(81, 14).Value = -32122
$ws.Cells.Item(84, 8).Value = 5781.8184
$ws.Cells.Item(84, 9).Value = 2325
$ws.Cells.Item(84, 10).Value = 15000
$ws.Cells.Item(84, 11).Value = 23250
$ws.Cells.Item(84, 12).Value = 150000
$ws.Cells.Item(84, 13).Value = -17946
$ws.Cells.Item(84, 14).Value = -160608
$ws.Cells.Item(107, 8).Value = 666.55554
$ws.Cells.Item(107, 10).Value = 809.6
$ws.Cells.Item(107, 12).Value = 2428.8
$ws.Cells.Item(107, 14).Value = -6268.8
$ws.Cells.Item(122, 8).Value = 3959.4517
$ws.Cells.Item(122, 9).Value = 2002.7894
$ws.Cells.Item(122, 10).Value = 7057.5
$ws.Cells.Item(122, 11).Value = 6008.3682
$ws.Cells.Item(122, 12).Value = 21172.5
$ws.Cells.Item(122, 13).Value = -3558.3682
$ws.Cells.Item(122, 14).Value = -26072.5
$ws.Cells.Item(132, 8).Value = 18001.666
$ws.Cells.Item(132, 9).Value = 10000
$ws.Cells.Item(132, 10).Value = 22002.5
$ws.Cells.Item(132, 11).Value = 30000
$ws.Cells.Item(132, 12).Value = 66007.5
$ws.Cells.Item(132, 13).Value = -27470
$ws.Cells.Item(132, 14).Value = -71067.5
$ws.Cells.Item(136, 8).Value = 3905.5
$ws.Cells.Item(136, 9).Value = 2085.5
$ws.Cells.Item(136, 11).Value = 6256.5
$ws.Cells.Item(136, 13).Value = -3706.5
